# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet tab was renamed from SCD0276 to SCD0017
$ws.Name = "SCD0017"

# TC_ID cells (B2/B3) were updated from the old JIRA id DGS-291
# to the new test-case id SCD0017-006
$ws.Range("B2").Value = "SCD0017-006"
$ws.Range("B3").Value = "SCD0017-006"

# Column B needs to widen to fit the new, longer TC_ID text
# (bestFit-style width recompute for the new value)
$ws.Columns("B").ColumnWidth = 11.451822916666666

# Cursor/selection ended up on B4 after the edit
$ws.Range("B4").Select()
